$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the activity name (was "Research about proposals")
$ws.Range("A2").Value = "chart and presentatipn"

# Update Date Start / Date End values (serials keep the existing date format
# and avoid picking up a time-of-day component)
$ws.Range("B2").Value = 43728
$ws.Range("C2").Value = 43733

# Move the active selection to B5, matching the saved workbook state
$ws.Range("B5").Select()
